# Insert a new data row at row 131, shifting the existing rows 131-172 down
# to 132-173, then populate the new row 131 with its data.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row before the current row 131 (pushes 131..172 -> 132..173)
$ws.Rows(131).Insert()

# Populate the freshly inserted row 131 with the new record.
$ws.Range("A131").Value = 3
$ws.Range("B131").Value = "Femacal de La Calera"
$ws.Range("C131").Value = "Coquimbo"
$ws.Range("D131").Value = 44508
$ws.Range("E131").Value = 5
$ws.Range("F131").Value = 100112001
$ws.Range("G131").Value = "Berenjena"
$ws.Range("H131").Value = "Sin especificar"
$ws.Range("I131").Value = "Primera"
$ws.Range("J131").Value = 95
$ws.Range("K131").Value = 7000
$ws.Range("L131").Value = 7500
$ws.Range("M131").Value = 7263
$ws.Range("N131").Value = "$/caja 60 unidades"
$ws.Range("O131").Value = "Región de Arica y Parinacota"
$ws.Range("P131").Value = 121
$ws.Range("Q131").Value = 60
$ws.Range("R131").Value = "Hortaliza"

# Match the date-format style used by the rest of column D.
$ws.Range("D131").NumberFormat = $ws.Range("D130").NumberFormat
